$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Summary" section header row (old row 38); this shifts all
# subsequent rows up by one (old row 44 -> new row 43).
$ws.Rows.Item(38).Delete()

# Re-label the per-branch detail rows so each line includes the branch name,
# e.g. "     New nominations" -> "     Civilian, New nominations".
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Confirmed "
$ws.Range("A9").Value  = "     Civilian, Unconfirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "

$ws.Range("A12").Value = "     Other Civilian, New nominations"
$ws.Range("A13").Value = "     Other Civilian, Confirmed "
$ws.Range("A14").Value = "     Other Civilian, Unconfirmed "

$ws.Range("A16").Value = "     Air Force, New nominations"
$ws.Range("A17").Value = "     Air Force, Confirmed "
$ws.Range("A18").Value = "     Air Force, Unconfirmed "

$ws.Range("A20").Value = "     Army, New nominations"
$ws.Range("A21").Value = "     Army, Confirmed "
$ws.Range("A22").Value = "     Army, Unconfirmed "
$ws.Range("A23").Value = "     Army, Withdrawn "

$ws.Range("A25").Value = "     Navy, New nominations"
$ws.Range("A26").Value = "     Navy, Confirmed "
$ws.Range("A27").Value = "     Navy, Unconfirmed "
$ws.Range("A28").Value = "     Navy, Withdrawn "

$ws.Range("A30").Value = "     Marine Corps, New nominations"
$ws.Range("A31").Value = "     Marine Corps, Confirmed "
$ws.Range("A32").Value = "     Marine Corps, Unconfirmed "

$ws.Range("A34").Value = "     Space Force, New nominations"
$ws.Range("A35").Value = "     Space Force, Confirmed "
$ws.Range("A36").Value = "     Space Force, Unconfirmed "
$ws.Range("A37").Value = "     Space Force, Withdrawn "

# The former "Total nominations received this Session" / "Total nominations
# carried over from the First Session" rows are renamed and swap order, and
# their number formats need to swap too (26351 uses thousands separators,
# 0 does not). Use copy/paste-special of formats from existing cells that
# already carry the desired style so we reuse the workbook's existing style
# records instead of minting new ones.
$ws.Range("B7").Copy()                 # B7 already uses the plain "General" style
$ws.Range("B39").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B12").Copy()                # B12 already uses the "#,##0" style
$ws.Range("B38").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A38").Value = "Total new nominations"
$ws.Range("B38").Value = 26351
$ws.Range("A39").Value = "Total carryover nominations"
$ws.Range("B39").Value = 0

# Rows 40-43 keep their existing labels/values (Total confirmed, Total
# unconfirmed, Total withdrawn, Total returned to the White House); nothing
# else to change there since they simply shifted up from 41-44.
